$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parameters")

# theta_step parameter value change: 50 -> 5
$ws.Range("C6").Value = 5

# reflect the author's final cursor position after the edit
$ws.Activate()
$null = $ws.Range("C9").Select()
